$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 273146.16
$ws.Range("I62").Value = 5381.9
$ws.Range("J62").Value = 451655.66
$ws.Range("K62").Value = 5381.9
$ws.Range("L62").Value = 451655.66
$ws.Range("M62").Value = -4757.9
$ws.Range("N62").Value = -452903.66
$ws.Range("H65").Value = 273146.16
$ws.Range("I65").Value = 5381.9
$ws.Range("J65").Value = 451655.66
$ws.Range("K65").Value = 26909.5
$ws.Range("L65").Value = 2258278.3
$ws.Range("M65").Value = -23789.5
$ws.Range("N65").Value = -2264518.3
$ws.Range("H74").Value = 4074.2
$ws.Range("I74").Value = 3468.3635
$ws.Range("K74").Value = 3468.3635
$ws.Range("M74").Value = -2532.3635
$ws.Range("H77").Value = 4074.2
$ws.Range("I77").Value = 3468.3635
$ws.Range("K77").Value = 17341.8175
$ws.Range("M77").Value = -12661.8175
$ws.Range("H107").Value = 1891.6666
$ws.Range("I107").Value = 337.83334
$ws.Range("J107").Value = 4999.3335
$ws.Range("K107").Value = 337.83334
$ws.Range("L107").Value = 4999.3335
$ws.Range("M107").Value = 1582.16666
$ws.Range("N107").Value = -8839.333500000001
$ws.Range("H109").Value = 81441.664
$ws.Range("J109").Value = 81441.664
$ws.Range("L109").Value = 81441.664
$ws.Range("N109").Value = -84215.664
$ws.Range("H118").Value = 572.25
$ws.Range("I118").Value = 572.25
$ws.Range("K118").Value = 1716.75
$ws.Range("M118").Value = -59.75
$ws.Range("H132").Value = 1815.8966
$ws.Range("I132").Value = 1347.12
$ws.Range("K132").Value = 4041.36
$ws.Range("M132").Value = -1511.36
$ws.Range("H133").Value = 93371.664
$ws.Range("J133").Value = 93371.664
$ws.Range("L133").Value = 93371.664
$ws.Range("N133").Value = -103491.664
$ws.Range("H134").Value = 54207.145
$ws.Range("J134").Value = 54207.145
$ws.Range("L134").Value = 54207.145
$ws.Range("N134").Value = -64347.145
$ws.Range("H135").Value = 1256.8182
$ws.Range("I135").Value = 1256.8182
$ws.Range("K135").Value = 11311.3638
$ws.Range("M135").Value = -8776.363799999999
$ws.Range("H136").Value = 96491.664
$ws.Range("J136").Value = 96491.664
$ws.Range("L136").Value = 96491.664
$ws.Range("N136").Value = -106691.664
$ws.Range("H139").Value = 99990
$ws.Range("J139").Value = 99990
$ws.Range("L139").Value = 99990
$ws.Range("N139").Value = -110270
$ws.Range("H140").Value = 84697.71000000001
$ws.Range("J140").Value = 83980.664
$ws.Range("L140").Value = 83980.664
$ws.Range("N140").Value = -94340.664

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 249.63637
$ws.Range("I4").Value = 286.5
$ws.Range("J4").Value = 151.33333
$ws.Range("K4").Value = 286.5
$ws.Range("L4").Value = 151.33333
$ws.Range("M4").Value = -170.5
$ws.Range("N4").Value = -383.33333
$ws.Range("H132").Value = 1877.8611
$ws.Range("I132").Value = 1643.9375
$ws.Range("K132").Value = 4931.8125
$ws.Range("M132").Value = -2401.8125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 80771.42999999999
$ws.Range("J132").Value = 80771.42999999999
$ws.Range("L132").Value = 80771.42999999999
$ws.Range("N132").Value = -90891.42999999999
$ws.Range("H135").Value = 48608.668
$ws.Range("J135").Value = 48608.668
$ws.Range("L135").Value = 48608.668
$ws.Range("N135").Value = -58748.668
$ws.Range("H138").Value = 95750.664
$ws.Range("J138").Value = 95750.664
$ws.Range("L138").Value = 95750.664
$ws.Range("N138").Value = -106030.664
$ws.Range("H140").Value = 124244.75
$ws.Range("J140").Value = 70565.42999999999
$ws.Range("L140").Value = 70565.42999999999
$ws.Range("N140").Value = -80925.42999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 188.48148
$ws.Range("J7").Value = 308.25
$ws.Range("L7").Value = 308.25
$ws.Range("N7").Value = -534.25
$ws.Range("H31").Value = 3976.923
$ws.Range("J31").Value = 5558.6
$ws.Range("L31").Value = 5558.6
$ws.Range("N31").Value = -6148.6
$ws.Range("H34").Value = 3976.923
$ws.Range("J34").Value = 5558.6
$ws.Range("L34").Value = 5558.6
$ws.Range("N34").Value = -5962.6
$ws.Range("H103").Value = 27666.666
$ws.Range("I103").Value = 1499.5
$ws.Range("K103").Value = 1499.5
$ws.Range("M103").Value = -327.5
$ws.Range("H132").Value = 1907.3529
$ws.Range("I132").Value = 1687.0834
$ws.Range("K132").Value = 5061.2502
$ws.Range("M132").Value = -2531.2502
$ws.Range("H141").Value = 208622.83
$ws.Range("J141").Value = 228513.33
$ws.Range("L141").Value = 228513.33
$ws.Range("N141").Value = -238873.33

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 812.6
$ws.Range("J34").Value = 1125.1428
$ws.Range("L34").Value = 3375.4284
$ws.Range("N34").Value = -3543.4284
$ws.Range("H68").Value = 994.5
$ws.Range("J68").Value = 990
$ws.Range("L68").Value = 2970
$ws.Range("N68").Value = -4592
$ws.Range("H71").Value = 994.5
$ws.Range("J71").Value = 990
$ws.Range("L71").Value = 8910
$ws.Range("N71").Value = -17022
$ws.Range("H122").Value = 1443872.4
$ws.Range("I122").Value = 999
$ws.Range("J122").Value = 2021021.8
$ws.Range("K122").Value = 8991
$ws.Range("L122").Value = 18189196.2
$ws.Range("M122").Value = -6541
$ws.Range("N122").Value = -18194096.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 2000
$ws.Range("I10").Value = 2000
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 2000
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()
$ws.Range("M10").Value = -1831
$ws.Range("H44").Value = 24258.666
$ws.Range("I44").Value = 7000
$ws.Range("J44").Value = 27710.4
$ws.Range("K44").Value = 7000
$ws.Range("L44").Value = 27710.4
$ws.Range("M44").Value = -6404
$ws.Range("N44").Value = -28902.4
$ws.Range("H52").Value = 25250
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 25250
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 25250
$ws.Range("M52").ClearContents()
$ws.Range("N52").Value = -25768
$ws.Range("H113").Value = 3335193.2
$ws.Range("J113").Value = 5557891
$ws.Range("L113").Value = 5557891
$ws.Range("N113").Value = -5562231

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 17999.666
$ws.Range("I3").Value = 14000
$ws.Range("J3").Value = 19999.5
$ws.Range("K3").Value = 14000
$ws.Range("L3").Value = 19999.5
$ws.Range("M3").Value = -13888
$ws.Range("N3").Value = -20223.5
$ws.Range("H15").Value = 17999.666
$ws.Range("I15").Value = 14000
$ws.Range("J15").Value = 19999.5
$ws.Range("K15").Value = 14000
$ws.Range("L15").Value = 19999.5
$ws.Range("M15").Value = -13830
$ws.Range("N15").Value = -20339.5
$ws.Range("H22").Value = 1008.5455
$ws.Range("I22").Value = 848
$ws.Range("J22").Value = 1100.2858
$ws.Range("K22").Value = 848
$ws.Range("L22").Value = 1100.2858
$ws.Range("M22").Value = -553
$ws.Range("N22").Value = -1690.2858
$ws.Range("H27").Value = 1008.5455
$ws.Range("I27").Value = 848
$ws.Range("J27").Value = 1100.2858
$ws.Range("K27").Value = 848
$ws.Range("L27").Value = 1100.2858
$ws.Range("M27").Value = -741
$ws.Range("N27").Value = -1314.2858
$ws.Range("H40").Value = 5054063.5
$ws.Range("I40").Value = 3738.3333
$ws.Range("K40").Value = 3738.3333
$ws.Range("M40").Value = -3602.3333
$ws.Range("H87").Value = 27879.4
$ws.Range("J87").Value = 18600
$ws.Range("L87").Value = 18600
$ws.Range("N87").Value = -20846
$ws.Range("H90").Value = 27879.4
$ws.Range("J90").Value = 18600
$ws.Range("L90").Value = 55800
$ws.Range("N90").Value = -67032
$ws.Range("H132").Value = 8616.478999999999
$ws.Range("I132").Value = 11355.134
$ws.Range("K132").Value = 34065.402
$ws.Range("M132").Value = -31535.402

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 28124.6
$ws.Range("I45").Value = 14000
$ws.Range("J45").Value = 31655.75
$ws.Range("K45").Value = 14000
$ws.Range("L45").Value = 31655.75
$ws.Range("M45").Value = -13509
$ws.Range("N45").Value = -32637.75
$ws.Range("H108").Value = 100000
$ws.Range("J108").Value = 100000
$ws.Range("L108").Value = 100000
$ws.Range("N108").Value = -107680
$ws.Range("H122").Value = 6442.125
$ws.Range("I122").Value = 6371.75
$ws.Range("K122").Value = 19115.25
$ws.Range("M122").Value = -16665.25
